$d = $word.ActiveDocument

# Update p-values table: ecological, socioeconomic, Cod, Hake rows
$d.Content.Find.Execute("0.93", $true, $false, $false, $false, $false, $true, 1, $false, "0.87", 2)
$d.Content.Find.Execute("0.45", $true, $false, $false, $false, $false, $true, 1, $false, "0.42", 2)
$d.Content.Find.Execute("0.35", $true, $false, $false, $false, $false, $true, 1, $false, "0.33", 2)
$d.Content.Find.Execute("0.14", $true, $false, $false, $false, $false, $true, 1, $false, "0.97", 2)
